$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (SpawnDistance), shifting it (and GrabGage) right
$ws.Columns.Item(8).Insert()

# New column H header and value
$ws.Cells.Item(1, 8).Value = "SkillDistance"
$ws.Cells.Item(2, 8).Value = 130

# Match final column widths (G:H = 14.75, I = 17, no bestFit)
$ws.Columns.Item(7).ColumnWidth = 14
$ws.Columns.Item(8).ColumnWidth = 14
$ws.Columns.Item(9).ColumnWidth = 16.3

[void]$ws.Range("I5").Select()

Write-Host "done"
